$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update recalculated "desestacionalizada" (and a few other) values for existing rows ---
$ws.Cells.Item(110, 3).Value = 88.26000000000001
$ws.Cells.Item(111, 3).Value = 87.36
$ws.Cells.Item(112, 3).Value = 86.87
$ws.Cells.Item(113, 3).Value = 86.48
$ws.Cells.Item(114, 3).Value = 87
$ws.Cells.Item(115, 3).Value = 88.93000000000001
$ws.Cells.Item(116, 3).Value = 85.01000000000001
$ws.Cells.Item(117, 3).Value = 89.79000000000001
$ws.Cells.Item(118, 3).Value = 90.26000000000001
$ws.Cells.Item(119, 3).Value = 90.81999999999999
$ws.Cells.Item(120, 3).Value = 88.88
$ws.Cells.Item(121, 3).Value = 88.23
$ws.Cells.Item(123, 3).Value = 89.34
$ws.Cells.Item(124, 3).Value = 82.76000000000001
$ws.Cells.Item(125, 3).Value = 90.63
$ws.Cells.Item(126, 3).Value = 88.65000000000001
$ws.Cells.Item(127, 3).Value = 91.43000000000001
$ws.Cells.Item(128, 3).Value = 92.13
$ws.Cells.Item(129, 3).Value = 93.69
$ws.Cells.Item(130, 3).Value = 90.84999999999999
$ws.Cells.Item(132, 3).Value = 91.73
$ws.Cells.Item(133, 3).Value = 91.95
$ws.Cells.Item(134, 3).Value = 95.06999999999999
$ws.Cells.Item(135, 3).Value = 90.42
$ws.Cells.Item(136, 3).Value = 93.48999999999999
$ws.Cells.Item(137, 3).Value = 93.20999999999999
$ws.Cells.Item(139, 3).Value = 92.36
$ws.Cells.Item(140, 3).Value = 87.66
$ws.Cells.Item(142, 3).Value = 94.56999999999999
$ws.Cells.Item(143, 3).Value = 94.36
$ws.Cells.Item(144, 3).Value = 94.31999999999999
$ws.Cells.Item(145, 3).Value = 97.31999999999999
$ws.Cells.Item(146, 3).Value = 95.36
$ws.Cells.Item(147, 3).Value = 94.25
$ws.Cells.Item(148, 3).Value = 96.09
$ws.Cells.Item(149, 3).Value = 97.37
$ws.Cells.Item(150, 3).Value = 96.81
$ws.Cells.Item(151, 3).Value = 95.91
$ws.Cells.Item(152, 3).Value = 93.42
$ws.Cells.Item(153, 3).Value = 97.84999999999999
$ws.Cells.Item(154, 3).Value = 97.26000000000001
$ws.Cells.Item(155, 3).Value = 97.53
$ws.Cells.Item(156, 3).Value = 97.95
$ws.Cells.Item(158, 3).Value = 101.25
$ws.Cells.Item(159, 3).Value = 99.88
$ws.Cells.Item(160, 3).Value = 100.68
$ws.Cells.Item(161, 3).Value = 97.11
$ws.Cells.Item(163, 3).Value = 99.05
$ws.Cells.Item(164, 3).Value = 101.49
$ws.Cells.Item(165, 3).Value = 101.02
$ws.Cells.Item(166, 3).Value = 99.88
$ws.Cells.Item(167, 3).Value = 100.62
$ws.Cells.Item(168, 3).Value = 101.59
$ws.Cells.Item(169, 3).Value = 99.41
$ws.Cells.Item(170, 3).Value = 99.91
$ws.Cells.Item(172, 3).Value = 99.41
$ws.Cells.Item(173, 3).Value = 101.18
$ws.Cells.Item(174, 3).Value = 101.96
$ws.Cells.Item(176, 3).Value = 99.64
$ws.Cells.Item(177, 3).Value = 98.8
$ws.Cells.Item(179, 3).Value = 99.81999999999999
$ws.Cells.Item(180, 3).Value = 98.69
$ws.Cells.Item(181, 3).Value = 100.26
$ws.Cells.Item(184, 3).Value = 99.14
$ws.Cells.Item(185, 3).Value = 101.38
$ws.Cells.Item(187, 3).Value = 103.91
$ws.Cells.Item(188, 3).Value = 99.97
$ws.Cells.Item(190, 3).Value = 100.8
$ws.Cells.Item(191, 3).Value = 101.1
$ws.Cells.Item(193, 3).Value = 99.44
$ws.Cells.Item(194, 3).Value = 99.81
$ws.Cells.Item(195, 3).Value = 100.03
$ws.Cells.Item(196, 3).Value = 102.55
$ws.Cells.Item(197, 3).Value = 98.65000000000001
$ws.Cells.Item(198, 3).Value = 99.40000000000001
$ws.Cells.Item(199, 3).Value = 99.31
$ws.Cells.Item(200, 3).Value = 98.22
$ws.Cells.Item(200, 5).Value = 98.39
$ws.Cells.Item(201, 3).Value = 98.23
$ws.Cells.Item(204, 3).Value = 99.68000000000001
$ws.Cells.Item(208, 3).Value = 91.23999999999999
$ws.Cells.Item(209, 3).Value = 97.33
$ws.Cells.Item(210, 3).Value = 97.69
$ws.Cells.Item(211, 3).Value = 98.03
$ws.Cells.Item(211, 5).Value = 94.26000000000001
$ws.Cells.Item(212, 3).Value = 100.07
$ws.Cells.Item(212, 5).Value = 96.54000000000001
$ws.Cells.Item(214, 3).Value = 100.84
$ws.Cells.Item(217, 3).Value = 100.1
$ws.Cells.Item(218, 3).Value = 101.21
$ws.Cells.Item(219, 3).Value = 103.96
$ws.Cells.Item(220, 3).Value = 102.35
$ws.Cells.Item(221, 3).Value = 101.4
$ws.Cells.Item(222, 5).Value = 101.38
$ws.Cells.Item(223, 3).Value = 103.13
$ws.Cells.Item(224, 3).Value = 101.96
$ws.Cells.Item(225, 3).Value = 101.34
$ws.Cells.Item(226, 3).Value = 100.5
$ws.Cells.Item(228, 3).Value = 104.59
$ws.Cells.Item(229, 3).Value = 103.03
$ws.Cells.Item(230, 3).Value = 100.54
$ws.Cells.Item(231, 3).Value = 99.77
$ws.Cells.Item(232, 3).Value = 100.23
$ws.Cells.Item(233, 3).Value = 102.68
$ws.Cells.Item(234, 3).Value = 102.88
$ws.Cells.Item(235, 3).Value = 101.12
$ws.Cells.Item(236, 3).Value = 101.58
$ws.Cells.Item(236, 5).Value = 101.4
$ws.Cells.Item(237, 3).Value = 103.08
$ws.Cells.Item(237, 5).Value = 97.67
$ws.Cells.Item(238, 3).Value = 101.22
$ws.Cells.Item(239, 3).Value = 98.31999999999999
$ws.Cells.Item(240, 3).Value = 101.42
$ws.Cells.Item(240, 5).Value = 100.93
$ws.Cells.Item(241, 3).Value = 104.84
$ws.Cells.Item(243, 3).Value = 101.69
$ws.Cells.Item(244, 3).Value = 101.07
$ws.Cells.Item(245, 3).Value = 98.12
$ws.Cells.Item(245, 5).Value = 94.48
$ws.Cells.Item(246, 3).Value = 98.31
$ws.Cells.Item(246, 5).Value = 93.59
$ws.Cells.Item(247, 3).Value = 97.34999999999999
$ws.Cells.Item(247, 5).Value = 87.59999999999999
$ws.Cells.Item(248, 3).Value = 97.97
$ws.Cells.Item(249, 3).Value = 98.41
$ws.Cells.Item(250, 3).Value = 100.44
$ws.Cells.Item(251, 3).Value = 102.2
$ws.Cells.Item(252, 3).Value = 100.51
$ws.Cells.Item(253, 3).Value = 100.49
$ws.Cells.Item(254, 3).Value = 101.84
$ws.Cells.Item(255, 3).Value = 101.61
$ws.Cells.Item(256, 3).Value = 103.02
$ws.Cells.Item(257, 3).Value = 103.23
$ws.Cells.Item(258, 3).Value = 101.51
$ws.Cells.Item(258, 5).Value = 102.94
$ws.Cells.Item(259, 2).Value = 100.56
$ws.Cells.Item(259, 3).Value = 104.36

# --- Append new row 260 for 01-07-2021 ---
$ws.Range("A260").NumberFormat = "@"
$ws.Range("A260").Value = "01-07-2021"
$ws.Range("A260").Style = "Normal"
$ws.Cells.Item(260, 2).Value = 102.63
$ws.Cells.Item(260, 3).Value = 103.38
$ws.Cells.Item(260, 5).Value = 106.08
